# Swap the data of row 5 and row 6 for the columns that differ between
# the two records (A, B, D, E, F, G, H, P, Q, R, Z, AB). The remaining
# columns (C, I, K, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY)
# are identical between the two rows, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 gets row 6's distinguishing values ---
$ws.Range("A5").Value = 112128524
$ws.Range("B5").Value = 90814
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P5").Value = "Godmyr (Godmyr), Ly lm"
$ws.Range("Q5").Value = 690281
$ws.Range("R5").Value = 7126404
$ws.Range("Z5").Value = "15:00"
$ws.Range("AB5").Value = "15:00"

# --- Row 6 gets row 5's (original) distinguishing values ---
$ws.Range("A6").Value = 112127546
$ws.Range("B6").Value = 90806
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 4361
$ws.Range("F6").Value = "Orange taggsvamp"
$ws.Range("G6").Value = "Hydnellum aurantiacum"
$ws.Range("H6").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P6").Value = "Svarvarmyran (Svarvarmyran), Ly lm"
$ws.Range("Q6").Value = 690408
$ws.Range("R6").Value = 7125570
$ws.Range("Z6").Value = "14:25"
$ws.Range("AB6").Value = "14:25"
